$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.963.34"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "1.906.46"
$ws.Range("E3").Value = "  -4.35%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").Value = "'324.93"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("E8").Value = "  -3.32%  "
$ws.Range("D9").Value = "'45.58"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").Value = "'0.07749"
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").Value = "'0.9828"
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("D12").Value = "'22.09"
$ws.Range("E12").Value = "  -3.55%  "
$ws.Range("D13").Value = "1.969.75"
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("D14").Value = "'6.999"
$ws.Range("E14").Value = "  -3.83%  "
$ws.Range("D15").Value = "'5.685"
$ws.Range("E15").Value = "  -3.10%  "
$ws.Range("D16").Value = "'0.07050"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "'84.21"
$ws.Range("E18").Value = "  -5.05%  "
$ws.Range("D19").Value = "'0.000009574"
$ws.Range("E19").Value = "  -4.08%  "
$ws.Range("E20").Value = "  -3.65%  "
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").Value = "28.929.07"
$ws.Range("E22").Value = "  -2.20%  "
$ws.Range("D23").Value = "'5.340"
$ws.Range("E23").Value = "  -3.50%  "
$ws.Range("E24").Value = "  -2.81%  "
$ws.Range("D25").Value = "2.137.85"
$ws.Range("E25").Value = "  -4.88%  "
$ws.Range("D26").Value = "'2.080"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").Value = "'156.72"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "'19.14"
$ws.Range("E28").Value = "  -2.67%  "
$ws.Range("D29").Value = "'5.607"
$ws.Range("E29").Value = "  -6.48%  "
$ws.Range("D30").Value = "'117.90"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").Value = "'1.833"
$ws.Range("E31").Value = "  -6.39%  "
$ws.Range("D32").Value = "'0.09276"
$ws.Range("D33").Value = "'0.8624"
$ws.Range("E33").Value = "  -4.96%  "
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("D35").Value = "'1.253"
$ws.Range("E35").Value = "  -7.21%  "
$ws.Range("D36").Value = "'3.018"
$ws.Range("E36").Value = "  -5.03%  "
$ws.Range("D37").Value = "'0.05711"
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("D38").Value = "'1.146"
$ws.Range("E38").Value = "  -2.24%  "
$ws.Range("D39").Value = "'1.003"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("E41").Value = "  -4.82%  "
$ws.Range("D42").Value = "'0.5538"
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("D44").Value = "'9.342"
$ws.Range("E44").Value = "  -4.93%  "
$ws.Range("D45").Value = "'2.737"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("D47").Value = "'11.35"
$ws.Range("E47").Value = "  -5.19%  "
$ws.Range("D48").Value = "'2.099"
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("D49").Value = "'0.06819"
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("D50").Value = "'111.62"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("D51").Value = "'0.000002591"
$ws.Range("E51").Value = "  -26.57%  "
